# Update points for phone 09876543 -> 120.00
#
# The sheet currently ends at row 35 with the "09876543" phone-number
# record (stored as text, to keep its leading zero) having 0 points.
# We keep that historical record as a plain numeric phone lookup key
# (09876543 -> 9876543) on row 35, and add a new row 36 that holds the
# original text-formatted phone number "09876543" together with its
# updated point total of 120.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 35: phone id becomes a genuine number (loses the leading zero),
# points stay at 0.
$ws.Range("A35").Value = 9876543
$ws.Range("C35").Value = 0

# Row 36 (new): re-add the phone number as text (apostrophe prefix keeps
# the leading zero / text formatting) with the updated point total.
$ws.Range("A36").Value = "'09876543"
$ws.Range("C36").Value = 120
